$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.349.81"
$ws.Range("E2").Value = "  +8.78%  "
$ws.Range("D3").Value = "'1.599.37"
$ws.Range("E3").Value = "  +8.13%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  -0.38%  "
$ws.Range("D5").Value = "'0.9968"
$ws.Range("E5").Value = "  +2.14%  "
$ws.Range("D6").Value = "'288.76"
$ws.Range("E6").Value = "  +3.30%  "
$ws.Range("D7").Value = "'0.3688"
$ws.Range("E7").Value = "  +0.83%  "
$ws.Range("D8").Value = "'0.3392"
$ws.Range("E8").Value = "  +10.31%  "
$ws.Range("D9").Value = "'42.57"
$ws.Range("E9").Value = "  +6.58%  "
$ws.Range("E10").Value = "  +7.18%  "
$ws.Range("D11").Value = "'0.07040"
$ws.Range("E11").Value = "  +5.85%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  -0.26%  "
$ws.Range("D13").Value = "'19.69"
$ws.Range("E13").Value = "  +8.54%  "
$ws.Range("E14").Value = "  +7.38%  "
$ws.Range("D15").Value = "'6.642"
$ws.Range("E15").Value = "  +7.20%  "
$ws.Range("D16").Value = "'0.9961"
$ws.Range("E16").Value = "  +1.90%  "
$ws.Range("D17").Value = "'1.597.16"
$ws.Range("E17").Value = "  +7.77%  "
$ws.Range("D18").Value = "'0.00001079"
$ws.Range("E18").Value = "  +4.90%  "
$ws.Range("D19").Value = "'0.06609"
$ws.Range("E19").Value = "  +11.40%  "
$ws.Range("D20").Value = "'78.07"
$ws.Range("E20").Value = "  +12.36%  "
$ws.Range("D21").Value = "'16.14"
$ws.Range("E21").Value = "  +11.00%  "
$ws.Range("D22").Value = "'6.026"
$ws.Range("E22").Value = "  +10.05%  "
$ws.Range("D23").Value = "'11.76"
$ws.Range("E23").Value = "  +6.54%  "
$ws.Range("D24").Value = "'22.337.85"
$ws.Range("E24").Value = "  +8.38%  "
$ws.Range("D25").Value = "'2.399"
$ws.Range("E25").Value = "  +6.57%  "
$ws.Range("D26").Value = "'2.496"
$ws.Range("E26").Value = "  +16.59%  "
$ws.Range("D27").Value = "'149.98"
$ws.Range("E27").Value = "  +6.25%  "
$ws.Range("D28").Value = "'19.53"
$ws.Range("E28").Value = "  +13.26%  "
$ws.Range("D29").Value = "'1.775.40"
$ws.Range("E29").Value = "  +8.33%  "
$ws.Range("D30").Value = "'120.34"
$ws.Range("E30").Value = "  +5.74%  "
$ws.Range("D31").Value = "'4.160"
$ws.Range("E31").Value = "  +4.47%  "
$ws.Range("D32").Value = "'6.014"
$ws.Range("E32").Value = "  +20.63%  "
$ws.Range("D33").Value = "'0.9464"
$ws.Range("E33").Value = "  +16.54%  "
$ws.Range("D34").Value = "'0.08242"
$ws.Range("E34").Value = "  +2.83%  "
$ws.Range("D35").Value = "'1.609"
$ws.Range("E35").Value = "  +4.06%  "
$ws.Range("D36").Value = "'5.306"
$ws.Range("E36").Value = "  +12.42%  "
$ws.Range("D37").Value = "'8.631"
$ws.Range("E37").Value = "  +11.07%  "
$ws.Range("D38").Value = "'11.75"
$ws.Range("E38").Value = "  +12.42%  "
$ws.Range("D39").Value = "'0.06154"
$ws.Range("E39").Value = "  +5.55%  "
$ws.Range("D40").Value = "'1.235"
$ws.Range("E40").Value = "  +0.70%  "
$ws.Range("D41").Value = "'0.02214"
$ws.Range("E41").Value = "  +8.13%  "
$ws.Range("D42").Value = "'0.2028"
$ws.Range("E42").Value = "  +7.58%  "
$ws.Range("D43").Value = "'0.9958"
$ws.Range("E43").Value = "  +1.98%  "
$ws.Range("D44").Value = "'0.5903"
$ws.Range("E44").Value = "  +11.48%  "
$ws.Range("D45").Value = "'13.14"
$ws.Range("E45").Value = "  +8.20%  "
$ws.Range("D46").Value = "'3.670"
$ws.Range("E46").Value = "  +4.00%  "
$ws.Range("D47").Value = "'0.5700"
$ws.Range("E47").Value = "  +9.67%  "
$ws.Range("D48").Value = "'125.71"
$ws.Range("E48").Value = "  +5.40%  "
$ws.Range("D49").Value = "'1.966"
$ws.Range("E49").Value = "  +9.25%  "
$ws.Range("D50").Value = "'0.06816"
$ws.Range("E50").Value = "  +5.46%  "
$ws.Range("D51").Value = "'73.59"
$ws.Range("E51").Value = "  +8.91%  "
